$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price/Volume) cells look numeric but are stored as text in the
# original workbook (t="inlineStr"). Force text storage so Excel does not
# auto-coerce the assigned string into a number, then restore "Normal" style
# so no residual number-format style is left attached to the cell.
$dRefs = @("D2", "D3", "D4", "D5", "D6", "D7", "D8", "D9", "D11", "D12", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D27", "D40", "D41", "D42", "D43", "D44", "D45", "D47", "D48", "D49")
foreach ($r in $dRefs) {
    $ws.Range($r).NumberFormat = "@"
}

$ws.Range("D2").Value = "237.23"
$ws.Range("D3").Value = "21.77"
$ws.Range("D4").Value = "5.380"
$ws.Range("D5").Value = "0.05623"
$ws.Range("D6").Value = "6.476"
$ws.Range("D7").Value = "3.349"
$ws.Range("D8").Value = "0.7969"
$ws.Range("D9").Value = "1.033"
$ws.Range("D11").Value = "0.07326"
$ws.Range("D12").Value = "0.03125"
$ws.Range("D13").Value = "0.02964"
$ws.Range("D14").Value = "0.09218"
$ws.Range("D15").Value = "0.001672"
$ws.Range("D16").Value = "3.262"
$ws.Range("D17").Value = "0.04771"
$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D18").Value = "0.0005709"
$ws.Range("E18").Value = "17OneONE"
$ws.Range("B19").Value = "TigerCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D19").Value = "0.006220"
$ws.Range("E19").Value = "18TigerCashTCH"
$ws.Range("B20").Value = "HotbitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D20").Value = "0.005078"
$ws.Range("E20").Value = "19HotbitTokenHTB"
$ws.Range("B21").Value = "BitKan"
$ws.Range("C21").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D21").Value = "0.001052"
$ws.Range("E21").Value = "20BitKanKAN"
$ws.Range("B22").Value = "NitroEx"
$ws.Range("C22").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D22").Value = "0.0001501"
$ws.Range("E22").Value = "21NitroExNTX"
$ws.Range("B23").Value = "UpBots"
$ws.Range("C23").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("D23").Value = "0.0003899"
$ws.Range("E23").Value = "22UpBotsUBXT"
$ws.Range("B24").Value = "LEO"
$ws.Range("C24").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D24").Value = "3.951"
$ws.Range("E24").Value = "23LEOLEOBestin24h"
$ws.Range("B25").Value = "BTSEToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D25").Value = "2.203"
$ws.Range("E25").Value = "24BTSETokenBTSE"
$ws.Range("D27").Value = "0.1061"
$ws.Range("D40").Value = "0.04082"
$ws.Range("D41").Value = "0.006974"
$ws.Range("D42").Value = "0.003003"
$ws.Range("D43").Value = "0.1040"
$ws.Range("D44").Value = "0.008830"
$ws.Range("D45").Value = "0.00005440"
$ws.Range("D47").Value = "0.6751"
$ws.Range("D48").Value = "0.03716"
$ws.Range("D49").Value = "0.00002100"

foreach ($r in $dRefs) {
    $ws.Range($r).Style = "Normal"
}
